$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, pushing existing rows 20-27 down to 21-28
$ws.Rows.Item(20).Insert()

# Fill in the new row 20 with data (same pattern as surrounding rows,
# new specific values for D, H, J, K, L, M, P)
$ws.Cells.Item(20, 1).Value = 4
$ws.Cells.Item(20, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(20, 3).Value = "Los Lagos"
$ws.Cells.Item(20, 4).Value = 44504
$ws.Cells.Item(20, 4).NumberFormat = $ws.Cells.Item(21, 4).NumberFormat
$ws.Cells.Item(20, 5).Value = 10
$ws.Cells.Item(20, 6).Value = 300000000
$ws.Cells.Item(20, 7).Value = "Espárragos"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 180
$ws.Cells.Item(20, 11).Value = 1600
$ws.Cells.Item(20, 12).Value = 1600
$ws.Cells.Item(20, 13).Value = 1600
$ws.Cells.Item(20, 14).Value = "`$/kilo"
$ws.Cells.Item(20, 15).Value = "Provincia de Linares"
$ws.Cells.Item(20, 16).Value = 1600
$ws.Cells.Item(20, 17).Value = 1
$ws.Cells.Item(20, 18).Value = "Hortaliza"
